# repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows to match re-pulled data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -10
    4  = -3
    5  = 8
    6  = -3
    7  = 3
    10 = -13
    11 = -4
    13 = -4
    20 = -5
    22 = 1
    23 = -6
    25 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
